$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.419.40'
$ws.Cells.Item(2, 5).Value = '  +1.61%  '
$ws.Cells.Item(3, 4).Value = '3.147.68'
$ws.Cells.Item(3, 5).Value = '  +1.00%  '
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).Value = '603.16'
$ws.Cells.Item(5, 5).Value = '  -0.64%  '
$ws.Cells.Item(6, 4).Value = '144.15'
$ws.Cells.Item(6, 5).Value = '  -0.23%  '
$ws.Cells.Item(7, 5).Value = '  -0.14%  '
$ws.Cells.Item(8, 4).Value = '3.139.94'
$ws.Cells.Item(8, 5).Value = '  +0.74%  '
$ws.Cells.Item(9, 4).Value = '0.524'
$ws.Cells.Item(9, 5).Value = '  +1.20%  '
$ws.Cells.Item(10, 4).Value = '0.150'
$ws.Cells.Item(10, 5).Value = '  +0.57%  '
$ws.Cells.Item(11, 5).Value = '  +3.52%  '
$ws.Cells.Item(12, 4).Value = '0.471'
$ws.Cells.Item(12, 5).Value = '  +1.00%  '
$ws.Cells.Item(13, 4).Value = '0.0000255'
$ws.Cells.Item(13, 5).Value = '  +2.48%  '
$ws.Cells.Item(14, 4).Value = '35.57'
$ws.Cells.Item(14, 5).Value = '  +1.34%  '
$ws.Cells.Item(15, 4).Value = '3.653.58'
$ws.Cells.Item(15, 5).Value = '  +0.72%  '
$ws.Cells.Item(16, 5).Value = '  +2.96%  '
$ws.Cells.Item(17, 4).Value = '64.413.16'
$ws.Cells.Item(17, 5).Value = '  +1.43%  '
$ws.Cells.Item(18, 4).Value = '3.143.41'
$ws.Cells.Item(18, 5).Value = '  +0.77%  '
$ws.Cells.Item(19, 4).Value = '6.90'
$ws.Cells.Item(19, 5).Value = '  +1.87%  '
$ws.Cells.Item(20, 4).Value = '482.53'
$ws.Cells.Item(20, 5).Value = '  +1.98%  '
$ws.Cells.Item(21, 4).Value = '14.63'
$ws.Cells.Item(21, 5).Value = '  +0.82%  '
$ws.Cells.Item(22, 4).Value = '0.712'
$ws.Cells.Item(22, 5).Value = '  +1.41%  '
$ws.Cells.Item(23, 4).Value = '7.69'
$ws.Cells.Item(23, 5).Value = '  -0.42%  '
$ws.Cells.Item(24, 4).Value = '88.07'
$ws.Cells.Item(24, 5).Value = '  +6.20%  '
$ws.Cells.Item(25, 4).Value = '13.45'
$ws.Cells.Item(25, 5).Value = '  +0.25%  '
$ws.Cells.Item(26, 4).Value = '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.05%  '
$ws.Cells.Item(27, 4).Value = '2.75'
$ws.Cells.Item(27, 5).Value = '  -0.71%  '
$ws.Cells.Item(28, 4).Value = '8.41'
$ws.Cells.Item(28, 5).Value = '  +0.35%  '
$ws.Cells.Item(29, 4).Value = '7.14'
$ws.Cells.Item(29, 5).Value = '  +4.67%  '
$ws.Cells.Item(30, 4).Value = '2.08'
$ws.Cells.Item(30, 5).Value = '  +0.58%  '
$ws.Cells.Item(31, 5).Value = '  -4.05%  '
$ws.Cells.Item(32, 5).Value = '  -0.15%  '
$ws.Cells.Item(33, 4).Value = '26.92'
$ws.Cells.Item(33, 5).Value = '  +3.16%  '
$ws.Cells.Item(34, 4).Value = '2.69'
$ws.Cells.Item(34, 5).Value = '  +1.25%  '
$ws.Cells.Item(35, 4).Value = '1.10'
$ws.Cells.Item(35, 5).Value = '  -1.14%  '
$ws.Cells.Item(36, 4).Value = '6.05'
$ws.Cells.Item(36, 5).Value = '  +2.47%  '
$ws.Cells.Item(37, 4).Value = '0.0₃0762'
$ws.Cells.Item(37, 5).Value = '  +1.21%  '
$ws.Cells.Item(38, 4).Value = '52.92'
$ws.Cells.Item(38, 5).Value = '  +0.39%  '
$ws.Cells.Item(39, 4).Value = '3.05'
$ws.Cells.Item(39, 5).Value = '  +4.46%  '
$ws.Cells.Item(40, 4).Value = '439.90'
$ws.Cells.Item(40, 5).Value = '  -2.89%  '
$ws.Cells.Item(41, 4).Value = '0.0396'
$ws.Cells.Item(41, 5).Value = '  +1.28%  '
$ws.Cells.Item(42, 5).Value = '  +1.43%  '
$ws.Cells.Item(43, 4).Value = '8.26'
$ws.Cells.Item(43, 5).Value = '  -0.36%  '
$ws.Cells.Item(44, 4).Value = '2.867.50'
$ws.Cells.Item(44, 5).Value = '  +1.19%  '
$ws.Cells.Item(45, 2).Value = 'TheGraph'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(45, 4).Value = '0.262'
$ws.Cells.Item(45, 5).Value = '  +0.05%  '
$ws.Cells.Item(46, 2).Value = 'Fetch.AI'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(46, 4).Value = '2.24'
$ws.Cells.Item(46, 5).Value = '  -1.28%  '
$ws.Cells.Item(47, 4).Value = '2.46'
$ws.Cells.Item(47, 5).Value = '  +3.13%  '
$ws.Cells.Item(48, 4).Value = '0.998'
$ws.Cells.Item(48, 5).Value = '  -0.05%  '
$ws.Cells.Item(49, 4).Value = '26.00'
$ws.Cells.Item(49, 5).Value = '  +0.01%  '
$ws.Cells.Item(50, 4).Value = '0.114'
$ws.Cells.Item(50, 5).Value = '  +1.02%  '
$ws.Cells.Item(51, 4).Value = '121.44'
$ws.Cells.Item(51, 5).Value = '  +2.28%  '
